# Updates odds/statistics values in Sheet1 to match the 2024-11-17 FlashScore
# refresh. Only numeric "Odd_*" cells changed; text/header columns (A-F) are
# untouched. Values below are grouped by row for readability.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.6
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 3.5
$ws.Range("L2").Value = 4
$ws.Range("Q2").Value = 2.88
$ws.Range("R2").Value = 1.4
$ws.Range("Z2").Value = 26
$ws.Range("AA2").Value = 26
$ws.Range("AH2").Value = 13
$ws.Range("AJ2").Value = 34
$ws.Range("AN2").Value = 4.33
$ws.Range("BA2").Value = 126
# Row 3
$ws.Range("G3").Value = 2.3
$ws.Range("I3").Value = 3.4
$ws.Range("AA3").Value = 23
$ws.Range("AG3").Value = 7.5
$ws.Range("AO3").Value = 15
$ws.Range("AU3").Value = 9
$ws.Range("AZ3").Value = 67
# Row 4
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("Q4").Value = 2.2
$ws.Range("R4").Value = 1.65
$ws.Range("AA4").Value = 19
$ws.Range("AG4").Value = 9
$ws.Range("AH4").Value = 17
$ws.Range("AN4").Value = 4
# Row 5
$ws.Range("L5").Value = 7
$ws.Range("N5").Value = 8.5
$ws.Range("R5").Value = 1.75
$ws.Range("S5").Value = 1.44
$ws.Range("T5").Value = 2.63
$ws.Range("AC5").Value = 8.5
$ws.Range("AI5").Value = 23
$ws.Range("AS5").Value = 201
$ws.Range("AT5").Value = 2.63
$ws.Range("AX5").Value = 41
$ws.Range("BA5").Value = 201
# Row 7
$ws.Range("G7").Value = 1.72
$ws.Range("H7").Value = 3.35
$ws.Range("I7").Value = 4.65
$ws.Range("J7").Value = 2.25
$ws.Range("K7").Value = 2.12
$ws.Range("O7").Value = 1.32
$ws.Range("P7").Value = 2.85
$ws.Range("Q7").Value = 1.93
$ws.Range("U7").Value = 1.83
$ws.Range("V7").Value = 1.78
$ws.Range("X7").Value = 7.8
$ws.Range("Y7").Value = 8
$ws.Range("Z7").Value = 13.5
$ws.Range("AA7").Value = 14.5
$ws.Range("AB7").Value = 28
$ws.Range("AC7").Value = 8.75
$ws.Range("AD7").Value = 6.6
$ws.Range("AE7").Value = 16
$ws.Range("AG7").Value = 11.5
$ws.Range("AH7").Value = 26
$ws.Range("AI7").Value = 15
$ws.Range("AJ7").Value = 80
$ws.Range("AN7").Value = 3.55
$ws.Range("AT7").Value = 2.72
$ws.Range("AW7").Value = 6.3
# Row 8
$ws.Range("G8").Value = 2.62
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 2.67
$ws.Range("J8").Value = 3.2
$ws.Range("K8").Value = 2.07
$ws.Range("L8").Value = 3.1
$ws.Range("N8").Value = 7.9
$ws.Range("O8").Value = 1.31
$ws.Range("P8").Value = 2.9
$ws.Range("Q8").Value = 1.98
$ws.Range("R8").Value = 1.75
$ws.Range("U8").Value = 1.65
$ws.Range("V8").Value = 1.98
$ws.Range("W8").Value = 8
$ws.Range("X8").Value = 13
$ws.Range("Y8").Value = 9.75
$ws.Range("Z8").Value = 30
$ws.Range("AA8").Value = 23
$ws.Range("AB8").Value = 32
$ws.Range("AC8").Value = 8.75
$ws.Range("AD8").Value = 5.8
$ws.Range("AE8").Value = 12.5
$ws.Range("AF8").Value = 55
$ws.Range("AG8").Value = 9
$ws.Range("AH8").Value = 14.5
$ws.Range("AI8").Value = 9.5
$ws.Range("AJ8").Value = 32
$ws.Range("AK8").Value = 22
$ws.Range("AL8").Value = 28
$ws.Range("AM8").Value = 400
$ws.Range("AN8").Value = 4.6
$ws.Range("AO8").Value = 14
$ws.Range("AP8").Value = 20
$ws.Range("AQ8").Value = 60
$ws.Range("AR8").Value = 90
$ws.Range("AS8").Value = 250
$ws.Range("AT8").Value = 2.6
$ws.Range("AU8").Value = 6.3
$ws.Range("AW8").Value = 4.7
$ws.Range("AX8").Value = 13.5
$ws.Range("AY8").Value = 18
$ws.Range("AZ8").Value = 55
$ws.Range("BA8").Value = 75
$ws.Range("BB8").Value = 200
# Row 9
$ws.Range("BD9").Value = 151
# Row 11
$ws.Range("G11").Value = 3.2
$ws.Range("I11").Value = 2.45
$ws.Range("J11").Value = 4
$ws.Range("L11").Value = 3.25
$ws.Range("M11").Value = 1.1
$ws.Range("N11").Value = 7
$ws.Range("X11").Value = 15
$ws.Range("AH11").Value = 10
$ws.Range("AI11").Value = 10
$ws.Range("AN11").Value = 5
# Row 12
$ws.Range("G12").Value = 3.7
$ws.Range("H12").Value = 3.5
$ws.Range("I12").Value = 1.95
$ws.Range("J12").Value = 4.75
$ws.Range("L12").Value = 2.63
$ws.Range("Q12").Value = 2.25
$ws.Range("R12").Value = 1.62
$ws.Range("W12").Value = 9
$ws.Range("X12").Value = 19
$ws.Range("AD12").Value = 7
$ws.Range("AG12").Value = 6
$ws.Range("AH12").Value = 8.5
$ws.Range("AI12").Value = 9
$ws.Range("AJ12").Value = 17
$ws.Range("AN12").Value = 6
$ws.Range("AO12").Value = 23
$ws.Range("AQ12").Value = 81
$ws.Range("AR12").Value = 126
$ws.Range("AU12").Value = 9
$ws.Range("AW12").Value = 3.75
$ws.Range("AX12").Value = 11
$ws.Range("AY12").Value = 23
# Row 13
$ws.Range("G13").Value = 3.7
$ws.Range("I13").Value = 2
$ws.Range("J13").Value = 4.33
$ws.Range("L13").Value = 2.75
$ws.Range("W13").Value = 9.5
$ws.Range("X13").Value = 17
$ws.Range("AQ13").Value = 67
# Row 15
$ws.Range("G15").Value = 1.67
$ws.Range("L15").Value = 6.5
$ws.Range("M15").Value = 1.11
$ws.Range("N15").Value = 6.5
$ws.Range("Z15").Value = 12
$ws.Range("AA15").Value = 17
$ws.Range("AG15").Value = 10
$ws.Range("AH15").Value = 26
$ws.Range("AJ15").Value = 67
# Row 18
$ws.Range("M18").Value = 1.08
$ws.Range("N18").Value = 8
$ws.Range("O18").Value = 1.4
$ws.Range("P18").Value = 2.75
# Row 19
$ws.Range("S19").Value = 1.5
$ws.Range("T19").Value = 2.5
$ws.Range("AP19").Value = 29
$ws.Range("AT19").Value = 2.5
# Row 20
$ws.Range("O20").Value = 1.44
$ws.Range("P20").Value = 2.63
